$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoopFilter LPF")

# Adjust the cutoff frequency input (intermediate decimation) from 50 to 25.
# All dependent formula cells (B4, B6, B8, B9, B10, B12, B13, B14, B20, B21, B22)
# recalculate automatically from this single input change.
$ws.Range("B2").Value = 25

$excel.CalculateFullRebuild()
